# Update "request list.xlsx" per commit "[updated] update request list"
#
# The active sheet is "0.9.0_20150626" (the 2nd sheet). The edit:
#  - row 7 (previously blank) gets new content describing NVR OEM log collection
#  - row 5's description is reworded
#  - column C is widened a bit to fit the new text
#  - the saved selection moves to F10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row: "NVR OEM버전 로그수집" / "NVR OEM버전에서도 로그 수집할 수 있는 기능"
# (set this first so new shared strings are appended in the same order as the
# authored workbook: NVR string before the reworded row 5 string)
$ws.Range("C7").Value = "NVR OEM버전 로그수집"

# Reword row 5's description
$ws.Range("D5").Value = "프로그램 실행 중일때도 로그파일 수집할 수 있는 기능"

$ws.Range("D7").Value = "NVR OEM버전에서도 로그 수집할 수 있는 기능"

# Widen column C to fit the longer header text
$ws.Columns.Item(3).ColumnWidth = 22.43

# Update the saved cursor/selection position
$ws.Range("F10").Select()
